$wb = $excel.ActiveWorkbook

# Work on the BDSBaPCF sheet (2nd sheet) - remove the "municipal solid waste"
# row's link to natural gas peaker (B9) and replace with a plain 0 value.
$ws = $wb.Worksheets.Item("BDSBaPCF")
$ws.Range("B17").Value = 0

# Make BDSBaPCF the active/selected sheet and set the new selection there.
$ws.Activate()
$ws.Range("B18").Select()
